$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update phone numbers to include the "602" Cali area code and fix malformed "+572" prefixes.
$ws.Range("E10").Value = "+57 602 3966729"
$ws.Range("E8").Value = "+57 602 6630509"
$ws.Range("E26").Value = "+57 602 3341781"
$ws.Range("E27").Value = "+57 602 3341781"
